$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.685.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.455.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.10"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.506"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.152"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.99%  "
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.343"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.82"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.904.15"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.532.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.454.33"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.20"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.43%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.89"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.07%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.38%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.582.90"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0845"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.62%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.35"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.07%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "432.35"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.90%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.69"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("B35").Value = "POPCAT"
$ws.Range("C35").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.01"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +103.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.83"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.02"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.110"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.98"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.307"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.54"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.10"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.10"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.62"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.38"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.491"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.94%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0718"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.565"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0917"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.01%  "
